# Update workbook to reflect carjacking data through 2022-07-16
# (adds the extra day of data collected on 2022-07-24).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab from "...07-15" to "...07-16"
$ws.Name = "Through 2022-07-16"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 07-16)"

# July 2022 year-to-date count increases from 82 to 91
$ws.Range("I8").Value = 91

# Full-year 2022 total increases from 887 to 896
$ws.Range("I14").Value = 896
